$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cells in existing rows 255-287 ---
$ws.Cells.Item(255, 4).Value = 44748
$ws.Cells.Item(255, 10).Value = 500
$ws.Cells.Item(255, 11).Value = 19000
$ws.Cells.Item(255, 12).Value = 20000
$ws.Cells.Item(255, 13).Value = 19500
$ws.Cells.Item(255, 16).Value = 1950
$ws.Cells.Item(256, 4).Value = 44246
$ws.Cells.Item(256, 10).Value = 600
$ws.Cells.Item(256, 11).Value = 14000
$ws.Cells.Item(256, 12).Value = 15000
$ws.Cells.Item(256, 13).Value = 14500
$ws.Cells.Item(256, 16).Value = 1450
$ws.Cells.Item(257, 4).Value = 44641
$ws.Cells.Item(257, 10).Value = 540
$ws.Cells.Item(257, 11).Value = 19000
$ws.Cells.Item(257, 12).Value = 20000
$ws.Cells.Item(257, 13).Value = 19500
$ws.Cells.Item(257, 16).Value = 1950
$ws.Cells.Item(258, 4).Value = 44699
$ws.Cells.Item(258, 10).Value = 440
$ws.Cells.Item(259, 4).Value = 44596
$ws.Cells.Item(259, 10).Value = 560
$ws.Cells.Item(259, 11).Value = 18000
$ws.Cells.Item(259, 12).Value = 19000
$ws.Cells.Item(259, 13).Value = 18500
$ws.Cells.Item(259, 16).Value = 1850
$ws.Cells.Item(260, 4).Value = 44291
$ws.Cells.Item(260, 10).Value = 520
$ws.Cells.Item(260, 11).Value = 12000
$ws.Cells.Item(260, 12).Value = 13000
$ws.Cells.Item(260, 13).Value = 12500
$ws.Cells.Item(260, 16).Value = 1250
$ws.Cells.Item(261, 4).Value = 44568
$ws.Cells.Item(261, 10).Value = 600
$ws.Cells.Item(261, 11).Value = 18000
$ws.Cells.Item(261, 12).Value = 19000
$ws.Cells.Item(261, 13).Value = 18500
$ws.Cells.Item(261, 16).Value = 1850
$ws.Cells.Item(262, 4).Value = 44620
$ws.Cells.Item(262, 10).Value = 500
$ws.Cells.Item(262, 11).Value = 19000
$ws.Cells.Item(262, 12).Value = 20000
$ws.Cells.Item(262, 13).Value = 19500
$ws.Cells.Item(262, 16).Value = 1950
$ws.Cells.Item(263, 4).Value = 44711
$ws.Cells.Item(263, 10).Value = 480
$ws.Cells.Item(263, 11).Value = 18000
$ws.Cells.Item(263, 12).Value = 19000
$ws.Cells.Item(263, 13).Value = 18500
$ws.Cells.Item(263, 16).Value = 1850
$ws.Cells.Item(264, 4).Value = 44239
$ws.Cells.Item(264, 11).Value = 14500
$ws.Cells.Item(264, 12).Value = 15000
$ws.Cells.Item(264, 13).Value = 14750
$ws.Cells.Item(264, 16).Value = 1475
$ws.Cells.Item(265, 4).Value = 44272
$ws.Cells.Item(265, 10).Value = 600
$ws.Cells.Item(265, 11).Value = 13000
$ws.Cells.Item(265, 12).Value = 14000
$ws.Cells.Item(265, 13).Value = 13500
$ws.Cells.Item(265, 16).Value = 1350
$ws.Cells.Item(266, 4).Value = 44348
$ws.Cells.Item(266, 10).Value = 400
$ws.Cells.Item(266, 11).Value = 12000
$ws.Cells.Item(266, 12).Value = 12500
$ws.Cells.Item(266, 13).Value = 12250
$ws.Cells.Item(266, 16).Value = 1225
$ws.Cells.Item(267, 4).Value = 44533
$ws.Cells.Item(267, 10).Value = 600
$ws.Cells.Item(267, 11).Value = 19000
$ws.Cells.Item(267, 12).Value = 20000
$ws.Cells.Item(267, 13).Value = 19500
$ws.Cells.Item(267, 16).Value = 1950
$ws.Cells.Item(268, 4).Value = 44229
$ws.Cells.Item(268, 10).Value = 400
$ws.Cells.Item(268, 11).Value = 12000
$ws.Cells.Item(268, 12).Value = 13000
$ws.Cells.Item(268, 13).Value = 12500
$ws.Cells.Item(268, 16).Value = 1250
$ws.Cells.Item(269, 4).Value = 44505
$ws.Cells.Item(269, 10).Value = 660
$ws.Cells.Item(270, 4).Value = 44746
$ws.Cells.Item(270, 10).Value = 480
$ws.Cells.Item(270, 11).Value = 19000
$ws.Cells.Item(270, 12).Value = 20000
$ws.Cells.Item(270, 13).Value = 19500
$ws.Cells.Item(270, 16).Value = 1950
$ws.Cells.Item(271, 4).Value = 44386
$ws.Cells.Item(271, 10).Value = 700
$ws.Cells.Item(271, 11).Value = 12000
$ws.Cells.Item(271, 12).Value = 12500
$ws.Cells.Item(271, 13).Value = 12250
$ws.Cells.Item(271, 16).Value = 1225
$ws.Cells.Item(272, 4).Value = 44690
$ws.Cells.Item(272, 10).Value = 500
$ws.Cells.Item(273, 4).Value = 44631
$ws.Cells.Item(273, 10).Value = 560
$ws.Cells.Item(273, 11).Value = 19000
$ws.Cells.Item(273, 12).Value = 20000
$ws.Cells.Item(273, 13).Value = 19500
$ws.Cells.Item(273, 16).Value = 1950
$ws.Cells.Item(274, 4).Value = 44235
$ws.Cells.Item(274, 10).Value = 500
$ws.Cells.Item(274, 11).Value = 16000
$ws.Cells.Item(274, 12).Value = 17000
$ws.Cells.Item(274, 13).Value = 16500
$ws.Cells.Item(274, 16).Value = 1650
$ws.Cells.Item(275, 4).Value = 44582
$ws.Cells.Item(275, 11).Value = 19000
$ws.Cells.Item(275, 12).Value = 20000
$ws.Cells.Item(275, 13).Value = 19500
$ws.Cells.Item(275, 16).Value = 1950
$ws.Cells.Item(276, 4).Value = 44307
$ws.Cells.Item(276, 10).Value = 600
$ws.Cells.Item(276, 11).Value = 11500
$ws.Cells.Item(276, 12).Value = 12000
$ws.Cells.Item(276, 13).Value = 11750
$ws.Cells.Item(276, 16).Value = 1175
$ws.Cells.Item(277, 4).Value = 44344
$ws.Cells.Item(277, 11).Value = 12000
$ws.Cells.Item(277, 12).Value = 12500
$ws.Cells.Item(277, 13).Value = 12250
$ws.Cells.Item(277, 16).Value = 1225
$ws.Cells.Item(278, 4).Value = 44433
$ws.Cells.Item(278, 11).Value = 16000
$ws.Cells.Item(278, 12).Value = 17000
$ws.Cells.Item(278, 13).Value = 16500
$ws.Cells.Item(278, 16).Value = 1650
$ws.Cells.Item(279, 4).Value = 44265
$ws.Cells.Item(279, 10).Value = 660
$ws.Cells.Item(279, 11).Value = 13000
$ws.Cells.Item(279, 12).Value = 14000
$ws.Cells.Item(279, 13).Value = 13500
$ws.Cells.Item(279, 16).Value = 1350
$ws.Cells.Item(280, 4).Value = 44747
$ws.Cells.Item(280, 10).Value = 400
$ws.Cells.Item(280, 11).Value = 19000
$ws.Cells.Item(280, 12).Value = 20000
$ws.Cells.Item(280, 13).Value = 19500
$ws.Cells.Item(280, 16).Value = 1950
$ws.Cells.Item(281, 4).Value = 44747
$ws.Cells.Item(281, 10).Value = 400
$ws.Cells.Item(281, 11).Value = 20000
$ws.Cells.Item(281, 12).Value = 21000
$ws.Cells.Item(281, 13).Value = 20500
$ws.Cells.Item(281, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(281, 16).Value = 2050
$ws.Cells.Item(282, 4).Value = 44421
$ws.Cells.Item(282, 10).Value = 720
$ws.Cells.Item(283, 4).Value = 44215
$ws.Cells.Item(283, 10).Value = 440
$ws.Cells.Item(283, 11).Value = 13000
$ws.Cells.Item(283, 12).Value = 14000
$ws.Cells.Item(283, 13).Value = 13500
$ws.Cells.Item(283, 16).Value = 1350
$ws.Cells.Item(284, 4).Value = 44566
$ws.Cells.Item(284, 10).Value = 560
$ws.Cells.Item(284, 11).Value = 18000
$ws.Cells.Item(284, 12).Value = 19000
$ws.Cells.Item(284, 13).Value = 18500
$ws.Cells.Item(284, 16).Value = 1850
$ws.Cells.Item(285, 4).Value = 44244
$ws.Cells.Item(285, 10).Value = 600
$ws.Cells.Item(285, 11).Value = 15000
$ws.Cells.Item(285, 12).Value = 16000
$ws.Cells.Item(285, 13).Value = 15500
$ws.Cells.Item(285, 16).Value = 1550
$ws.Cells.Item(286, 4).Value = 44442
$ws.Cells.Item(286, 11).Value = 16000
$ws.Cells.Item(286, 12).Value = 17000
$ws.Cells.Item(286, 13).Value = 16500
$ws.Cells.Item(286, 16).Value = 1650
$ws.Cells.Item(287, 4).Value = 44663
$ws.Cells.Item(287, 10).Value = 460
$ws.Cells.Item(287, 11).Value = 19000
$ws.Cells.Item(287, 12).Value = 20000
$ws.Cells.Item(287, 13).Value = 19500
$ws.Cells.Item(287, 16).Value = 1950
# --- Append new rows 288-290 ---
$ws.Cells.Item(288, 1).Value = 8
$ws.Cells.Item(288, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(288, 3).Value = "Coquimbo"
$ws.Cells.Item(288, 4).Value = 44187
$ws.Cells.Item(288, 5).Value = 4
$ws.Cells.Item(288, 6).Value = 100112003
$ws.Cells.Item(288, 7).Value = "Ajo"
$ws.Cells.Item(288, 8).Value = "Chino"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 560
$ws.Cells.Item(288, 11).Value = 11000
$ws.Cells.Item(288, 12).Value = 12000
$ws.Cells.Item(288, 13).Value = 11500
$ws.Cells.Item(288, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(288, 15).Value = "China"
$ws.Cells.Item(288, 16).Value = 1150
$ws.Cells.Item(288, 17).Value = 10
$ws.Cells.Item(288, 18).Value = "Hortaliza"
$ws.Cells.Item(288, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(289, 1).Value = 8
$ws.Cells.Item(289, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(289, 3).Value = "Coquimbo"
$ws.Cells.Item(289, 4).Value = 44519
$ws.Cells.Item(289, 5).Value = 4
$ws.Cells.Item(289, 6).Value = 100112003
$ws.Cells.Item(289, 7).Value = "Ajo"
$ws.Cells.Item(289, 8).Value = "Chino"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 600
$ws.Cells.Item(289, 11).Value = 19000
$ws.Cells.Item(289, 12).Value = 20000
$ws.Cells.Item(289, 13).Value = 19500
$ws.Cells.Item(289, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(289, 15).Value = "China"
$ws.Cells.Item(289, 16).Value = 1950
$ws.Cells.Item(289, 17).Value = 10
$ws.Cells.Item(289, 18).Value = "Hortaliza"
$ws.Cells.Item(289, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(290, 1).Value = 8
$ws.Cells.Item(290, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(290, 3).Value = "Coquimbo"
$ws.Cells.Item(290, 4).Value = 44194
$ws.Cells.Item(290, 5).Value = 4
$ws.Cells.Item(290, 6).Value = 100112003
$ws.Cells.Item(290, 7).Value = "Ajo"
$ws.Cells.Item(290, 8).Value = "Chino"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 600
$ws.Cells.Item(290, 11).Value = 12000
$ws.Cells.Item(290, 12).Value = 13000
$ws.Cells.Item(290, 13).Value = 12500
$ws.Cells.Item(290, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(290, 15).Value = "China"
$ws.Cells.Item(290, 16).Value = 1250
$ws.Cells.Item(290, 17).Value = 10
$ws.Cells.Item(290, 18).Value = "Hortaliza"
$ws.Cells.Item(290, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
